# Update quarterly balance figures per "Actualizacion desde MV -datos-"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = 171
$ws.Range("C17").Value = 164
$ws.Range("C20").Value = 195
$ws.Range("C21").Value = 199
$ws.Range("C24").Value = 238
$ws.Range("C27").Value = 307
$ws.Range("C28").Value = 317
$ws.Range("C29").Value = 337
$ws.Range("C32").Value = 421
$ws.Range("I32").Value = 1088
$ws.Range("C35").Value = 466
$ws.Range("C36").Value = 504
$ws.Range("I36").Value = 1505
$ws.Range("C39").Value = 536
$ws.Range("C40").Value = 557
$ws.Range("I41").Value = 1940
$ws.Range("C43").Value = 664
$ws.Range("C44").Value = 640
$ws.Range("C45").Value = 695
$ws.Range("C46").Value = 715
$ws.Range("C47").Value = 762
$ws.Range("C48").Value = 775
$ws.Range("C49").Value = 779
$ws.Range("C50").Value = 786
$ws.Range("C52").Value = 891
$ws.Range("C53").Value = 892
$ws.Range("C54").Value = 883
$ws.Range("C55").Value = 914
$ws.Range("C56").Value = 900
$ws.Range("C57").Value = 887
$ws.Range("C58").Value = 892
$ws.Range("C60").Value = 934
$ws.Range("C61").Value = 1064
$ws.Range("C62").Value = 999
$ws.Range("C63").Value = 995
$ws.Range("C64").Value = 1009
$ws.Range("C65").Value = 1006
$ws.Range("C66").Value = 1055
$ws.Range("C69").Value = 1114
$ws.Range("C72").Value = 1211
$ws.Range("C73").Value = 1219
$ws.Range("C74").Value = 1247
$ws.Range("C76").Value = 1274
$ws.Range("C77").Value = 1302
$ws.Range("C79").Value = 1361
$ws.Range("I79").Value = 6434
$ws.Range("C80").Value = 1440
$ws.Range("C81").Value = 1491
$ws.Range("C83").Value = 1586
$ws.Range("C84").Value = 1701
$ws.Range("C86").Value = 1789
$ws.Range("I86").Value = 7849
$ws.Range("C88").Value = 1883
$ws.Range("C89").Value = 1957
$ws.Range("C90").Value = 2031
$ws.Range("C91").Value = 2121
$ws.Range("C92").Value = 2119
$ws.Range("I93").Value = 10384
$ws.Range("C94").Value = 2295
$ws.Range("C95").Value = 2382
$ws.Range("C96").Value = 2424
$ws.Range("I96").Value = 11001
$ws.Range("C97").Value = 2498
$ws.Range("C98").Value = 2523
$ws.Range("I98").Value = 11217
$ws.Range("C99").Value = 2608
$ws.Range("I99").Value = 12043
$ws.Range("C100").Value = 2698
$ws.Range("C101").Value = 2761
$ws.Range("I101").Value = 13874
$ws.Range("C103").Value = 3070
$ws.Range("C104").Value = 3117
$ws.Range("I104").Value = 16151
$ws.Range("C105").Value = 3184
$ws.Range("I105").Value = 16382
$ws.Range("C106").Value = 3239
$ws.Range("C107").Value = 3316
$ws.Range("I107").Value = 17189
$ws.Range("C108").Value = 3511
$ws.Range("I108").Value = 17800
$ws.Range("C110").Value = 3729
$ws.Range("C111").Value = 3922
$ws.Range("I111").Value = 18248
$ws.Range("I112").Value = 18733
$ws.Range("C113").Value = 4188
$ws.Range("C114").Value = 4307
$ws.Range("I114").Value = 19813
$ws.Range("C115").Value = 4445
$ws.Range("C116").Value = 4619
$ws.Range("C117").Value = 4745
$ws.Range("I117").Value = 21817
$ws.Range("C118").Value = 4808
$ws.Range("I118").Value = 22641
$ws.Range("C119").Value = 4895
$ws.Range("C120").Value = 4967
$ws.Range("C121").Value = 5114
$ws.Range("I121").Value = 24857
$ws.Range("C122").Value = 5281
$ws.Range("I122").Value = 25547
$ws.Range("C123").Value = 5490
$ws.Range("I124").Value = 27930
$ws.Range("C125").Value = 5623
$ws.Range("I125").Value = 28296
$ws.Range("C126").Value = 5732
$ws.Range("C127").Value = 5740
$ws.Range("I127").Value = 28017
$ws.Range("C128").Value = 5912
$ws.Range("I128").Value = 28345
$ws.Range("I129").Value = 29140
$ws.Range("C130").Value = 6031
$ws.Range("I130").Value = 29738
$ws.Range("C131").Value = 6065
$ws.Range("I131").Value = 30692
$ws.Range("C132").Value = 6086
$ws.Range("I132").Value = 31395
$ws.Range("C133").Value = 6179
$ws.Range("I133").Value = 32394
$ws.Range("C134").Value = 6337
$ws.Range("I134").Value = 33586
$ws.Range("C135").Value = 6329
$ws.Range("I135").Value = 34222
$ws.Range("C136").Value = 6364
$ws.Range("I136").Value = 34853
$ws.Range("C137").Value = 6448
$ws.Range("I137").Value = 35508
$ws.Range("C138").Value = 6465
$ws.Range("I138").Value = 36736
$ws.Range("C139").Value = 6607
$ws.Range("I139").Value = 37180
$ws.Range("C140").Value = 6900
$ws.Range("I140").Value = 39801
$ws.Range("C141").Value = 7373
$ws.Range("I141").Value = 42631
$ws.Range("I142").Value = 47031
$ws.Range("C143").Value = 8793
$ws.Range("I143").Value = 51746
$ws.Range("C144").Value = 11603
$ws.Range("I144").Value = 65228
$ws.Range("C145").Value = 11875
$ws.Range("I145").Value = 67359
$ws.Range("C146").Value = 12782
$ws.Range("I146").Value = 71079
$ws.Range("C147").Value = 14632
$ws.Range("I147").Value = 78902
$ws.Range("AH147").Value = 3699
$ws.Range("AJ147").Value = 1366
$ws.Range("AK147").Value = 282493
